$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.240084767341614
$ws.Range("B1").Value = 1.031797051429749
$ws.Range("C1").Value = 3.10218358039856
$ws.Range("D1").Value = 3.213366508483887
$ws.Range("E1").Value = 0.9376339912414551
